$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new fuel-group rows ("Fossil Gases", "Fossil Liquids") into each of the
# three year blocks (2030 / 2040 / 2050), right after "Biogenic Gases" / "Biogenic Liquids" ---
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(31).Insert()
$ws.Rows.Item(34).Insert()

# --- Populate the newly inserted blank rows ---
# Row 7: Fossil Gases (2030)
$ws.Range("A7").Value = "Fossil Gases"
$ws.Range("B7").Value = [double]"2030"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = [double]"0.0005753006392829026"
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = [double]"5.54567252492964e-05"
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""

# Row 10: Fossil Liquids (2030)
$ws.Range("A10").Value = "Fossil Liquids"
$ws.Range("B10").Value = [double]"2030"
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = [double]"0.0268422477686697"
$ws.Range("G10").Value = [double]"0.0001113324693366"
$ws.Range("H10").Value = [double]"0.0067299634888682"
$ws.Range("I10").Value = [double]"0.0165521734152812"
$ws.Range("J10").Value = [double]"0.0004082898153107"
$ws.Range("K10").Value = [double]"0.0120421763211259"

# Row 19: Fossil Gases (2040)
$ws.Range("A19").Value = "Fossil Gases"
$ws.Range("B19").Value = [double]"2040"
$ws.Range("C19").Value = ""
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = [double]"0.0003096724302784673"
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = [double]"5.915836569318424e-05"
$ws.Range("J19").Value = ""
$ws.Range("K19").Value = ""

# Row 22: Fossil Liquids (2040)
$ws.Range("A22").Value = "Fossil Liquids"
$ws.Range("B22").Value = [double]"2040"
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = [double]"0.007325348721515"
$ws.Range("G22").Value = [double]"0.000119688594305"
$ws.Range("H22").Value = [double]"0.006349186668053"
$ws.Range("I22").Value = [double]"0.007474822992653"
$ws.Range("J22").Value = [double]"0.0003623206834181"
$ws.Range("K22").Value = [double]"0.0116305293046029"

# Row 31: Fossil Gases (2050)
$ws.Range("A31").Value = "Fossil Gases"
$ws.Range("B31").Value = [double]"2050"
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = [double]"1.816557443737047e-05"
$ws.Range("G31").Value = ""
$ws.Range("H31").Value = ""
$ws.Range("I31").Value = [double]"2.158302253128147e-05"
$ws.Range("J31").Value = ""
$ws.Range("K31").Value = ""

# Row 34: Fossil Liquids (2050)
$ws.Range("A34").Value = "Fossil Liquids"
$ws.Range("B34").Value = [double]"2050"
$ws.Range("C34").Value = ""
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = ""
$ws.Range("F34").Value = [double]"0.0004892418401882201"
$ws.Range("G34").Value = [double]"0.0001079540432089"
$ws.Range("H34").Value = [double]"0.0057176896227644"
$ws.Range("I34").Value = [double]"0.0013237224787817"
$ws.Range("J34").Value = [double]"0.0003118678534556"
$ws.Range("K34").Value = [double]"0.010784457730426"

# --- Update the handful of pre-existing cells whose figures changed ---
# (Biogenic Liquids "Pass Aviation" + every "Overall Demand" total row, which now also
#  folds in the new Fossil Gases / Fossil Liquids / Aviation figures)
$ws.Range("H9").Value = [double]"0.0007394122816492519"
$ws.Range("F13").Value = [double]"0.03030065694914777"
$ws.Range("G13").Value = [double]"0.0001271406625798332"
$ws.Range("H13").Value = [double]"0.007469376098141392"
$ws.Range("I13").Value = [double]"0.0192962524219816"
$ws.Range("J13").Value = [double]"0.0004757422642456062"
$ws.Range("K13").Value = [double]"0.01329400971635877"
$ws.Range("H21").Value = [double]"0.0009093343726798519"
$ws.Range("F25").Value = [double]"0.009567611865459695"
$ws.Range("G25").Value = [double]"0.000145458557245453"
$ws.Range("H25").Value = [double]"0.007258548466475288"
$ws.Range("I25").Value = [double]"0.009399480220949516"
$ws.Range("J25").Value = [double]"0.0004442644614519767"
$ws.Range("K25").Value = [double]"0.013061973368652"
$ws.Range("H32").Value = [double]"3.428120915942979e-11"
$ws.Range("H33").Value = [double]"0.001202549676586484"
$ws.Range("F37").Value = [double]"0.001635014250071598"
$ws.Range("G37").Value = [double]"0.0001540004434633084"
$ws.Range("H37").Value = [double]"0.006920285818056424"
$ws.Range("I37").Value = [double]"0.00197072177551618"
$ws.Range("J37").Value = [double]"0.0004172605287836186"
$ws.Range("K37").Value = [double]"0.01281649969459145"
